# Scheduled runner update: refresh computed profit figures on the
# per-job profit sheets (market price / profit recalculation pass).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 280319.12
$ws.Range("J17").Value = 280319.12
$ws.Range("L17").Value = 840957.36
$ws.Range("N17").Value = -841293.36

$ws.Range("H86").Value = 111116400
$ws.Range("I86").Value = 3521.2
$ws.Range("K86").Value = 3521.2
$ws.Range("M86").Value = -2398.2

$ws.Range("H89").Value = 111116400
$ws.Range("I89").Value = 3521.2
$ws.Range("K89").Value = 17606
$ws.Range("M89").Value = -11990

$ws.Range("H99").Value = 234.85715
$ws.Range("I99").Value = 224
$ws.Range("J99").Value = 300
$ws.Range("K99").Value = 672
$ws.Range("L99").Value = 900
$ws.Range("M99").Value = 826
$ws.Range("N99").Value = -3896

$ws.Range("H112").Value = 1529.5238
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1595.7894
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 4787.3682
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -7003.3682

$ws.Range("H123").Value = 25895
$ws.Range("J123").Value = 25895
$ws.Range("L123").Value = 25895
$ws.Range("N123").Value = -35695

$ws.Range("H127").Value = 2742.0256
$ws.Range("I127").Value = 1575
$ws.Range("J127").Value = 2875.4
$ws.Range("K127").Value = 4725
$ws.Range("L127").Value = 8626.200000000001
$ws.Range("M127").Value = 235
$ws.Range("N127").Value = -18546.2

$ws.Range("H138").Value = 3042.51
$ws.Range("I138").Value = 1209.7
$ws.Range("J138").Value = 3246.1555
$ws.Range("K138").Value = 3629.1
$ws.Range("L138").Value = 9738.466499999999
$ws.Range("M138").Value = 1510.9
$ws.Range("N138").Value = -20018.4665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9341.459000000001
$ws.Range("I32").Value = 8868.387000000001
$ws.Range("K32").Value = 8868.387000000001
$ws.Range("M32").Value = -8581.387000000001

$ws.Range("H45").Value = 2653.4
$ws.Range("I45").Value = 2392.6667
$ws.Range("K45").Value = 2392.6667
$ws.Range("M45").Value = -2015.6667

$ws.Range("H74").Value = 753.1842
$ws.Range("I74").Value = 771.05
$ws.Range("J74").Value = 733.3333
$ws.Range("K74").Value = 771.05
$ws.Range("L74").Value = 733.3333
$ws.Range("M74").Value = 102.95
$ws.Range("N74").Value = -2481.3333

$ws.Range("H77").Value = 753.1842
$ws.Range("I77").Value = 771.05
$ws.Range("J77").Value = 733.3333
$ws.Range("K77").Value = 3855.25
$ws.Range("L77").Value = 3666.6665
$ws.Range("M77").Value = 512.75
$ws.Range("N77").Value = -12402.6665

$ws.Range("H132").Value = 3563.238
$ws.Range("I132").Value = 3270.6155
$ws.Range("J132").Value = 4038.75
$ws.Range("K132").Value = 9811.8465
$ws.Range("L132").Value = 12116.25
$ws.Range("M132").Value = -7281.8465
$ws.Range("N132").Value = -17176.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 77112.14999999999
$ws.Range("I134").Value = 3707.6
$ws.Range("K134").Value = 11122.8
$ws.Range("M134").Value = -8587.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4045.3684
$ws.Range("I31").Value = 1178.3636
$ws.Range("J31").Value = 7987.5
$ws.Range("K31").Value = 1178.3636
$ws.Range("L31").Value = 7987.5
$ws.Range("M31").Value = -883.3635999999999
$ws.Range("N31").Value = -8577.5

$ws.Range("H34").Value = 4045.3684
$ws.Range("I34").Value = 1178.3636
$ws.Range("J34").Value = 7987.5
$ws.Range("K34").Value = 1178.3636
$ws.Range("L34").Value = 7987.5
$ws.Range("M34").Value = -976.3635999999999
$ws.Range("N34").Value = -8391.5

$ws.Range("H140").Value = 52020
$ws.Range("J140").Value = 52020
$ws.Range("L140").Value = 52020
$ws.Range("N140").Value = -62380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2358.3572
$ws.Range("I117").Value = 1064
$ws.Range("J117").Value = 2876.1
$ws.Range("K117").Value = 3192
$ws.Range("L117").Value = 8628.299999999999
$ws.Range("M117").Value = 250
$ws.Range("N117").Value = -15512.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3686.25
$ws.Range("I80").Value = 3632
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 3632
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -2634
$ws.Range("N80").Value = -6496

$ws.Range("H83").Value = 3686.25
$ws.Range("I83").Value = 3632
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 18160
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -13168
$ws.Range("N83").Value = -32484

$ws.Range("H113").Value = 1743.7333
$ws.Range("I113").Value = 1816.1111
$ws.Range("J113").Value = 1635.1666
$ws.Range("K113").Value = 1816.1111
$ws.Range("L113").Value = 1635.1666
$ws.Range("M113").Value = 353.8888999999999
$ws.Range("N113").Value = -5975.1666

$ws.Range("H132").Value = 4029.1428
$ws.Range("I132").Value = 2942
$ws.Range("J132").Value = 4633.1113
$ws.Range("K132").Value = 8826
$ws.Range("L132").Value = 13899.3339
$ws.Range("M132").Value = -6296
$ws.Range("N132").Value = -18959.3339

$ws.Range("H138").Value = 30494.285
$ws.Range("J138").Value = 30494.285
$ws.Range("L138").Value = 30494.285
$ws.Range("N138").Value = -40774.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7847.143
$ws.Range("I122").Value = 8659.166999999999
$ws.Range("J122").Value = 2975
$ws.Range("K122").Value = 25977.501
$ws.Range("L122").Value = 8925
$ws.Range("M122").Value = -23527.501
$ws.Range("N122").Value = -13825

$ws.Range("H132").Value = 2483.3845
$ws.Range("I132").Value = 2239.4482
$ws.Range("K132").Value = 6718.344599999999
$ws.Range("M132").Value = -4188.344599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 607.12
$ws.Range("I126").Value = 543.2222
$ws.Range("J126").Value = 771.4286
$ws.Range("K126").Value = 1629.6666
$ws.Range("L126").Value = 2314.2858
$ws.Range("M126").Value = 840.3334
$ws.Range("N126").Value = -7254.2858

$ws.Range("H132").Value = 1471.5172
$ws.Range("I132").Value = 986.34784
$ws.Range("J132").Value = 3331.3333
$ws.Range("K132").Value = 2959.04352
$ws.Range("L132").Value = 9993.999899999999
$ws.Range("M132").Value = -429.0435200000002
$ws.Range("N132").Value = -15053.9999
